# Fruta / hortaliza, semanal
# Insert a new weekly record at row 32 (shifting the existing rows 32-39 down to 33-40)
# and populate it with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32, pushing rows 32-39 down to 33-40.
$ws.Rows.Item(32).Insert()

# Fill in the new row 32 with this week's record.
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C32").Value = 'Ñuble'
$ws.Range("D32").Value = 44609
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 'Fruta'
$ws.Range("G32").Value = 100103
$ws.Range("H32").Value = 'Frutos de hueso (carozo)'
$ws.Range("I32").Value = 100103002
$ws.Range("J32").Value = 'Ciruela'
$ws.Range("K32").Value = 'Black Amber'
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 10000
$ws.Range("O32").Value = 11000
$ws.Range("P32").Value = 10500
$ws.Range("Q32").Value = '$/bandeja 18 kilos granel'
$ws.Range("R32").Value = 'Región de O''Higgins'
$ws.Range("S32").Value = 583
$ws.Range("T32").Value = 18
